$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('C2').Value = 'MSG: None

MSG: The decision process concluded without a definitive choice for Friday''s movie, resulting in no decision being made.
'
$ws.Range('C3').Value = 'MSG: None

MSG: The decision has been recorded as no movie selection.
'
$ws.Range('C4').Value = 'MSG: None

MSG: The decision on which movie to show on Friday has ended without a clear choice being made.
'
$ws.Range('D4').Value = 'no_decision, '
$ws.Range('C5').Value = 'MSG: None

MSG: The movie "Oppenheimer" has been successfully selected for acquisition.
'
$ws.Range('C6').Value = 'MSG: None

MSG: The decision has been recorded, and the rights for "Barbie" have been acquired.
'
$ws.Range('C7').Value = 'MSG: None

MSG: The decision resulted in no agreement about which movie to show on Friday.
'
$ws.Range('C8').Value = 'MSG: None

MSG: The decision has been recorded to acquire the rights for "Barbie."
'
$ws.Range('C9').Value = 'MSG: None

MSG: The decision has been recorded to acquire the rights for "Barbie." 
'
$ws.Range('C10').Value = 'MSG: None

MSG: The decision has been recorded, and the rights for the movie "Barbie" have been acquired.
'
$ws.Range('C11').Value = 'MSG: None

MSG: The decision has been recorded, indicating that there was no consensus on which movie to show on Friday.
'
$ws.Range('C12').Value = 'MSG: None

MSG: The decision has been recorded as no decision being made regarding the movie for Friday.
'
$ws.Range('C13').Value = 'MSG: None

MSG: The decision has been made to acquire the rights to "Barbie" for the screening on Friday.
'
$ws.Range('D13').Value = 'Barbie_was_selected, '
$ws.Range('C14').Value = 'MSG: None

MSG: The decision has been recorded as "no decision" about the movie to be shown on Friday.
'
$ws.Range('D14').Value = 'no_decision, '
$ws.Range('C15').Value = 'MSG: None

MSG: The conversation ended without a decision about which movie to play on Friday. Thus, the outcome is that no choice has been made.
'
$ws.Range('C16').Value = 'MSG: None

MSG: I have made the decision to call the no_decision function based on the criteria provided.
'
$ws.Range('C17').Value = 'MSG: None

MSG: The decision to acquire the rights for "Barbie" has been recorded.
'
$ws.Range('C18').Value = 'MSG: None

MSG: It seems there was no consensus reached about the movie to show on Friday, so a decision cannot be made. Therefore, I have recorded the outcome as a no decision situation.
'
$ws.Range('D18').Value = 'no_decision, '
$ws.Range('C19').Value = 'MSG: None

MSG: The decision about which movie to show on Friday was not reached, so I have confirmed that no decision will be made.
'
$ws.Range('D19').Value = 'no_decision, '
$ws.Range('C20').Value = 'MSG: None

MSG: The decision has been recorded to acquire the rights for "Barbie."
'
$ws.Range('C21').Value = 'MSG: None

MSG: The decision has been recorded, indicating that no movie was ultimately selected for Friday.
'
$ws.Range('C22').Value = 'MSG: None

MSG: The decision regarding Friday''s movie has not been made.
'
$ws.Range('C23').Value = 'MSG: None

MSG: The decision-making process ended without a clear agreement on which movie to show on Friday, so I will not be acquiring any movie rights.
'
$ws.Range('D23').Value = 'no_decision, '
$ws.Range('C24').Value = 'MSG: None

MSG: No decision was made regarding which movie to show on Friday.
'
$ws.Range('D24').Value = 'no_decision, '
$ws.Range('C25').Value = 'MSG: None

MSG: The decision has been made to acquire the rights for both movies, "Oppenheimer" and "Barbie."
'
$ws.Range('C26').Value = 'MSG: None

MSG: The decision to show a movie on Friday was not reached in the discussion, so the outcome is a no decision.
'
$ws.Range('C27').Value = 'MSG: None

MSG: The decision has been recorded. The selected movie for Friday is "Barbie."
'
$ws.Range('C28').Value = 'MSG: None

MSG: I have recorded the decision to acquire the rights for both movies "Barbie" and "Oppenheimer" as they will be shown in their entirety on Friday.
'
$ws.Range('C29').Value = 'MSG: None

MSG: The decision has been recorded, and "Barbie" will be the movie shown on Friday.
'
$ws.Range('C31').Value = 'MSG: None

MSG: The decision has been recorded as "no decision."
'
$ws.Range('C32').Value = 'MSG: None

MSG: The decision about which movie to show on Friday has resulted in no agreement.
'
$ws.Range('C33').Value = 'MSG: None

MSG: The decision has been recorded, and the rights to "Barbie" have been acquired.
'
$ws.Range('C34').Value = 'MSG: None

MSG: The decision-making committee did not reach a consensus on the movie to be shown on Friday, leading to no decision being made.
'
$ws.Range('C35').Value = 'MSG: None

MSG: The decision-making process concluded without arriving at an agreement on which movie to show on Friday.
'
$ws.Range('C36').Value = 'MSG: None

MSG: The decision has been recorded as no decision was reached regarding the selection of a movie.
'
$ws.Range('D36').Value = 'no_decision, '
$ws.Range('C37').Value = 'MSG: None

MSG: The decision has been made to acquire the rights for "Oppenheimer".
'
$ws.Range('C38').Value = 'MSG: None

MSG: The decision has been recorded, and no movie will be selected for Friday.
'
$ws.Range('C39').Value = 'MSG: None

MSG: The decision has been recorded as no decision.
'
$ws.Range('C40').Value = 'MSG: None

MSG: I have acquired the rights for both movies.
'
$ws.Range('C41').Value = 'MSG: None

MSG: The decision regarding the movie to be shown on Friday remains unresolved, and no specific choice has been made.
'
$ws.Range('C42').Value = 'MSG: None

MSG: The decision process has concluded without a selection for Friday''s movie.
'
$ws.Range('C43').Value = 'MSG: None

MSG: The decision has been recorded as no decision being made about Friday''s movie.
'
$ws.Range('C44').Value = 'MSG: None

MSG: The decision has been made to acquire the rights for "Barbie" to be shown on Friday.
'
$ws.Range('D44').Value = 'Barbie_was_selected, '
$ws.Range('C45').Value = 'MSG: None

MSG: The decision about what movie to show on Friday was not finalized, so no acquisition of rights will be taking place.
'
$ws.Range('C46').Value = 'MSG: None

MSG: I have recorded the decision to acquire the rights to both movies.
'
$ws.Range('C47').Value = 'MSG: None

MSG: The decision to acquire the rights for "Barbie" has been successfully recorded.
'
$ws.Range('C48').Value = 'MSG: None

MSG: The decision process concluded without reaching an agreement on a movie for Friday.
'
$ws.Range('C49').Value = 'MSG: None

MSG: The decision has been recorded as a no decision.
'
$ws.Range('C50').Value = 'MSG: None

MSG: The decision has been recorded to acquire the rights for "Oppenheimer."
'
$ws.Range('C51').Value = 'MSG: None

MSG: The decision to acquire the rights for both movies, "Barbie" and "Oppenheimer," has been confirmed.
'
$ws.Range('C52').Value = 'MSG: None

MSG: The decision regarding Friday''s movie ended without a definitive conclusion.
'
$ws.Range('C53').Value = 'MSG: None

MSG: The decision has been recorded as "no decision" regarding the choice of a movie for Friday.
'
$ws.Range('C54').Value = 'MSG: None

MSG: The decision to show a movie on Friday was not reached.
'
$ws.Range('C55').Value = 'MSG: None

MSG: I have recorded the decision to acquire the rights for "Barbie" as the movie for the assembly.
'
$ws.Range('C56').Value = 'MSG: None

MSG: The committee did not reach a decision about which movie to show on Friday, so I will call the no_decision function.
'
$ws.Range('C57').Value = 'MSG: None

MSG: The committee did not reach a consensus on which movie to show on Friday, so I have recorded no decision on the movie selection.
'
$ws.Range('C58').Value = 'MSG: None

MSG: The decision has been recorded as no decision was made regarding the movie to show on Friday.
'

Write-Host "Applied changes to" 56 "rows"